$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.499.73'
$ws.Range("E2").Value = '  +2.01%  '

$ws.Range("D3").Value = '1.856.44'
$ws.Range("E3").Value = '  +1.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.04'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6963'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07693'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.62'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.86%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07787'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.163'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.65%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.849.24'
$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '91.28'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6939'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.301'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.19%  '

$ws.Range("D17").Value = '29.493.31'
$ws.Range("E17").Value = '  +2.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008340'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.65%  '

$ws.Range("D19").Value = '2.106.41'
$ws.Range("E19").Value = '  +1.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '238.29'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.75'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9994'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.632'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.08%  '

$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("E25").Value = '  +1.45%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.89'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.890'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.27'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.38%  '

$ws.Range("E29").Value = '  -0.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.248'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.152'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.202'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05100'
$ws.Range("D33").ClearFormats()

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7740'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.887'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.28%  '

$ws.Range("E36").Value = '  +0.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.687'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").Value = '1.315.66'
$ws.Range("E38").Value = '  +7.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01877'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.54%  '

$ws.Range("E40").Value = '  +0.97%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9520'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.16'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.785'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.73%  '

$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.844'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.40%  '

$ws.Range("E46").Value = '  +2.57%  '

$ws.Range("D47").Value = '2.005.81'
$ws.Range("E47").Value = '  +1.48%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5218'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.787'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.17'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.959'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.93%  '
